$wb = $excel.ActiveWorkbook

# --- Collections sheet: add a "Link" column (after Comments, before Count Code) ---
# and a "Lifestage" column at the end, matching the Distributions sheet layout.
$ws = $wb.Worksheets.Item("Collections")

# Insert a new column H ("Link"); it inherits formatting/width info from the
# shift, same as the existing columns around it.
$ws.Columns.Item(8).Insert()
$ws.Cells.Item(3, 8).Value = "Link"
$ws.Columns.Item(8).ColumnWidth = $ws.Columns.Item(7).ColumnWidth

# Add the "Lifestage" column after Weight(g) (col M). Insert-then-delete the
# following column so the new cell correctly inherits the bordered header
# style used by the rest of row 3, without permanently shifting the columns
# beyond it.
$ws.Columns.Item(13).Insert()
$ws.Cells.Item(3, 13).Value = "Lifestage"
$ws.Columns.Item(14).Delete()

# --- View state: Collections becomes the active tab/selection, Distributions
# reverts to an unselected, scrolled-to-top view. ---
$ws.Activate()
$ws.Range("E9").Select()
